# Password Recovery System moves from "Not Started" (column A) to "Done" (column C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Not Started) currently lists tasks in rows 4-12, with A4 holding
# "Password Recovery System". Removing it shifts the remaining tasks up by one
# row and leaves the last row (12) blank.
for ($r = 4; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r + 1, 1).Value()
}
$ws.Cells.Item(12, 1).Value = ""

# Column C (Done) already has entries in rows 4-12; append the task as the
# next completed item in row 13.
$ws.Cells.Item(13, 3).Value = "Password Recovery System"

# Reflect the new selection/cursor position recorded in the workbook.
[void]$ws.Range("B12").Select()
